# daily auto push: 2025-10-08 06:41 UTC
# Append the day's new log row (row 79) to the bottom of the sei2 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

# Column A holds a date-like label but is stored as literal text (matching
# every prior row in the sheet, e.g. A2:A78), not a real Excel date serial.
# Assigning the plain string lets Excel's input parser auto-convert it to a
# date value, so we lead with an apostrophe to force text entry and then
# reset the cell style back to Normal (clearing the quote-prefix formatting
# Excel applied) so the cell ends up as plain text with no special style,
# exactly like the existing rows.
$ws.Range("A$row").Value = "'2025/10/08"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "水"
$ws.Range("C$row").Value = 14
$ws.Range("D$row").Value = 201
